# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 16-47: Tipo Doc, N Doc Trabajador, Nombre Trabajador, Periodo Mora, Valor Mora, Salario Basico
$rows = @(
    @(16, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1908", 20979, 828116),
    @(17, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1907", 33125, 828116),
    @(18, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1906", 33125, 828116),
    @(19, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1905", 33125, 828116),
    @(20, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1904", 33125, 828116),
    @(21, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1903", 33125, 828116),
    @(22, "CC", "1033375418", "GEIDER MANUEL PEREZ OYOLA", "1902", 33125, 828116),
    @(23, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1908", 20979, 828116),
    @(24, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1907", 33125, 828116),
    @(25, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1906", 33125, 828116),
    @(26, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1905", 33125, 828116),
    @(27, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1904", 33125, 828116),
    @(28, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1903", 33125, 828116),
    @(29, "CC", "1052735732", "GARIBALDIS GUERRERO FLOREZ", "1902", 33125, 828116),
    @(30, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1908", 20979, 828116),
    @(31, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1907", 33125, 828116),
    @(32, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1906", 33125, 828116),
    @(33, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1905", 33125, 828116),
    @(34, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1904", 33125, 828116),
    @(35, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1903", 33125, 828116),
    @(36, "CC", "1052731447", "JOSE MANUEL ALVAREZ DIAZ", "1902", 33125, 828116),
    @(37, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1908", 139205, 5494911),
    @(38, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1907", 219797, 5494911),
    @(39, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1906", 219797, 5494911),
    @(40, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1905", 219797, 5494911),
    @(41, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1904", 219797, 5494911),
    @(42, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1903", 219797, 5494911),
    @(43, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1902", 219797, 5494911),
    @(44, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1812", 219797, 5494911),
    @(45, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1811", 219797, 5494911),
    @(46, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1810", 219797, 5494911),
    @(47, "CC", "52501544", "PAOLA TATIANA OLIVER PEÑARANDA", "1808", 94889, 5494911)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[1]
    $ws.Range("C$rowNum").Value = $r[2]
    $ws.Range("D$rowNum").Value = $r[3]
    $ws.Range("E$rowNum").Value = $r[4]
    $ws.Range("F$rowNum").Value = $r[5]
    $ws.Range("G$rowNum").Value = $r[6]
}
